# Vragen.docx - "Querie update + meer cleanup"
#
# 1. Drop the stray "_GoBack" bookmark that sat (empty) right before the
#    "Titles / Principals / Crew" heading.
# 2. Fix the "flim" -> "film" typo in the map-question paragraph and drop
#    the now-stale spell-check proofErr markers around it.
# 3. Re-drop a fresh "_GoBack" bookmark at the very end of the document,
#    which is where Word leaves it after the cursor's last edit position.
#
# (Deleting the old "_GoBack" bookmark automatically renumbers the
#  remaining "_Hlk503442581" bookmark down from id 2 to id 1, and the
#  newly-added "_GoBack" bookmark at the end picks up id 2 again.)

$d = $word.ActiveDocument

# --- 1. Remove the obsolete "_GoBack" bookmark above "Titles / Principals / Crew"
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Correct "flim" -> "film" and let the replace swallow the old
#        spellStart/spellEnd proofErr markers that flagged the typo.
$found = $d.Content.Find.Execute("een flim speelt", $true, $false, $false, `
    $false, $false, $true, 1, $false, "een film speelt", 2)

# --- 3. Re-add "_GoBack" at the very end of the document (last paragraph).
$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)
$endRange = $lastParagraph.Range
$endRange.Collapse(0)
$endRange.Bookmarks.Add("_GoBack")
